# Rewrites the TANKERS section of the briefing:
#   - Track AR-YE (North): reorders SHELL/TEXACO/ARCO 2-1 entries, renames fields
#     (adds refueling type, splits callsign/altitude/freq layout, adds TKR tags),
#     and relocates the mid-paragraph _GoBack bookmark into the ARCO 2-1 line.
#   - Track AR-XC (South): same treatment for the 3-1 entries, and the trailing
#     bookmark-only paragraph becomes a plain blank paragraph.

$d = $word.ActiveDocument

# Locate the first paragraph of the block to replace ("- TEXACO 2-1 [KC-135] 57Y...")
# and the last paragraph of the block (the bookmark-only paragraph right before
# "Carrier S3-B").
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($null -eq $startPara -and $t.StartsWith("- TEXACO 2-1 [KC-135]")) {
        $startPara = $p
    }
    if ($t.StartsWith("Carrier S3-B")) {
        $endPara = $d.Paragraphs.Item($i - 1)
        break
    }
}

if ($null -eq $startPara -or $null -eq $endPara) {
    throw "Could not locate the TANKERS block paragraphs to replace."
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = '<w:p><w:r><w:t>- SHELL 2-1 [KC-135MPRS-drogue] 42</w:t></w:r><w:r><w:t>X</w:t></w:r><w:r><w:t>, FL220, 317.7</w:t></w:r><w:r><w:t>25</w:t></w:r><w:r><w:t xml:space="preserve"> AM</w:t></w:r></w:p><w:p><w:r><w:t>- TEXACO 2-1 [KC-135</w:t></w:r><w:r><w:t>-boom</w:t></w:r><w:r><w:t>] 5</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t>X</w:t></w:r><w:r><w:t>, FL240, 317.650 AM (TKR1</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>F16s)</w:t></w:r></w:p><w:p><w:r><w:t>- ARCO 2-1 [KC-130</w:t></w:r><w:r><w:t>-drogue</w:t></w:r><w:r><w:t xml:space="preserve">] </w:t></w:r><w:r><w:t>62</w:t></w:r><w:r><w:t>X</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">, FL200, </w:t></w:r><w:r><w:t>276.100</w:t></w:r><w:r><w:t xml:space="preserve"> AM</w:t></w:r><w:r><w:t xml:space="preserve"> (TKR2)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Track AR-XC (South)</w:t></w:r></w:p><w:p><w:r><w:t>- SHELL 3-1 [KC-135MPRS-drogue] 120</w:t></w:r><w:r><w:t>X</w:t></w:r><w:r><w:t xml:space="preserve">, FL220, </w:t></w:r><w:r><w:t>317.775</w:t></w:r><w:r><w:t xml:space="preserve"> AM</w:t></w:r><w:r><w:t xml:space="preserve"> (TKR1)</w:t></w:r></w:p><w:p><w:r><w:t>- TEXACO 3-1 [KC-135</w:t></w:r><w:r><w:t>-boom] 124</w:t></w:r><w:r><w:t>X</w:t></w:r><w:r><w:t>, FL240, 317.750 AM (TKR2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>F16s)</w:t></w:r></w:p><w:p><w:r><w:t>- ARCO 3-1 [KC-130</w:t></w:r><w:r><w:t>-drogue] 116</w:t></w:r><w:r><w:t>X</w:t></w:r><w:r><w:t>, FL200, 276.125 AM</w:t></w:r></w:p><w:p/>'

$rng.InsertXML($xml)
